$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Ensure columns B:E are formatted as Text so Excel does not reinterpret
# values like "25.917.37" or "1.000" as numbers/dates.
$ws.Range("B2:E51").NumberFormat = "@"

# Row 2
$ws.Cells.Item(2, 4).Value = "25.917.37"
$ws.Cells.Item(2, 5).Value = "  +0.15%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "1.733.24"
$ws.Cells.Item(3, 5).Value = "  -0.42%  "

# Row 4
$ws.Cells.Item(4, 4).Value = "0.9995"
$ws.Cells.Item(4, 5).Value = "  +0.00%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "245.74"
$ws.Cells.Item(5, 5).Value = "  +3.13%  "

# Row 6
$ws.Cells.Item(6, 4).Value = "0.9997"
$ws.Cells.Item(6, 5).Value = "  +0.00%  "

# Row 7
$ws.Cells.Item(7, 4).Value = "0.5035"
$ws.Cells.Item(7, 5).Value = "  -2.39%  "

# Row 8
$ws.Cells.Item(8, 4).Value = "0.2724"
$ws.Cells.Item(8, 5).Value = "  -0.43%  "

# Row 9
$ws.Cells.Item(9, 4).Value = "0.06173"
$ws.Cells.Item(9, 5).Value = "  +0.54%  "

# Row 10
$ws.Cells.Item(10, 4).Value = "1.737.71"
$ws.Cells.Item(10, 5).Value = "  +0.12%  "

# Row 11
$ws.Cells.Item(11, 4).Value = "0.07239"
$ws.Cells.Item(11, 5).Value = "  +1.01%  "

# Row 12
$ws.Cells.Item(12, 4).Value = "0.6539"
$ws.Cells.Item(12, 5).Value = "  +1.43%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "15.18"
$ws.Cells.Item(13, 5).Value = "  +1.67%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "4.782"
$ws.Cells.Item(14, 5).Value = "  +4.18%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "77.12"
$ws.Cells.Item(15, 5).Value = "  -0.20%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "1.0000"
$ws.Cells.Item(16, 5).Value = "  +0.00%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "0.9993"
$ws.Cells.Item(17, 5).Value = "  -0.01%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "25.914.52"
$ws.Cells.Item(18, 5).Value = "  +0.08%  "

# Row 19
$ws.Cells.Item(19, 4).Value = "11.87"
$ws.Cells.Item(19, 5).Value = "  +1.36%  "

# Row 20
$ws.Cells.Item(20, 4).Value = "0.000006817"
$ws.Cells.Item(20, 5).Value = "  +0.69%  "

# Row 21
$ws.Cells.Item(21, 2).Value = "Uniswap"
$ws.Cells.Item(21, 3).Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Cells.Item(21, 4).Value = "4.594"
$ws.Cells.Item(21, 5).Value = "  +7.70%  "

# Row 22
$ws.Cells.Item(22, 2).Value = "WrappedliquidstakedEther2.0"
$ws.Cells.Item(22, 3).Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Cells.Item(22, 4).Value = "1.963.23"
$ws.Cells.Item(22, 5).Value = "  +0.06%  "

# Row 23
$ws.Cells.Item(23, 4).Value = "8.792"
$ws.Cells.Item(23, 5).Value = "  +1.46%  "

# Row 24
$ws.Cells.Item(24, 4).Value = "5.477"
$ws.Cells.Item(24, 5).Value = "  +4.66%  "

# Row 25
$ws.Cells.Item(25, 4).Value = "133.93"
$ws.Cells.Item(25, 5).Value = "  -3.19%  "

# Row 26
$ws.Cells.Item(26, 5).Value = "  +0.84%  "

# Row 27
$ws.Cells.Item(27, 4).Value = "1.426"
$ws.Cells.Item(27, 5).Value = "  -5.44%  "

# Row 28
$ws.Cells.Item(28, 4).Value = "1.792"
$ws.Cells.Item(28, 5).Value = "  +1.72%  "

# Row 29
$ws.Cells.Item(29, 4).Value = "105.11"
$ws.Cells.Item(29, 5).Value = "  -0.66%  "

# Row 30
$ws.Cells.Item(30, 4).Value = "3.994"
$ws.Cells.Item(30, 5).Value = "  +1.04%  "

# Row 31
$ws.Cells.Item(31, 4).Value = "0.08121"
$ws.Cells.Item(31, 5).Value = "  -2.07%  "

# Row 32
$ws.Cells.Item(32, 5).Value = "  +1.60%  "

# Row 33
$ws.Cells.Item(33, 4).Value = "0.04723"
$ws.Cells.Item(33, 5).Value = "  +3.05%  "

# Row 34
$ws.Cells.Item(34, 4).Value = "2.655"
$ws.Cells.Item(34, 5).Value = "  -0.28%  "

# Row 35
$ws.Cells.Item(35, 4).Value = "0.9978"
$ws.Cells.Item(35, 5).Value = "  +0.94%  "

# Row 36
$ws.Cells.Item(36, 4).Value = "0.6131"
$ws.Cells.Item(36, 5).Value = "  -0.80%  "

# Row 37
$ws.Cells.Item(37, 4).Value = "2.743"
$ws.Cells.Item(37, 5).Value = "  +2.20%  "

# Row 38
$ws.Cells.Item(38, 4).Value = "0.8880"
$ws.Cells.Item(38, 5).Value = "  +20.25%  "

# Row 39
$ws.Cells.Item(39, 4).Value = "0.01604"
$ws.Cells.Item(39, 5).Value = "  -0.46%  "

# Row 40
$ws.Cells.Item(40, 4).Value = "1.959"
$ws.Cells.Item(40, 5).Value = "  +1.52%  "

# Row 41
$ws.Cells.Item(41, 4).Value = "0.9996"
$ws.Cells.Item(41, 5).Value = "  +0.03%  "

# Row 42
$ws.Cells.Item(42, 4).Value = "100.84"
$ws.Cells.Item(42, 5).Value = "  +3.24%  "

# Row 43
$ws.Cells.Item(43, 4).Value = "0.3923"
$ws.Cells.Item(43, 5).Value = "  +2.42%  "

# Row 44
$ws.Cells.Item(44, 4).Value = "5.022"
$ws.Cells.Item(44, 5).Value = "  +1.16%  "

# Row 45
$ws.Cells.Item(45, 4).Value = "0.1178"
$ws.Cells.Item(45, 5).Value = "  +4.72%  "

# Row 46
$ws.Cells.Item(46, 4).Value = "6.341"
$ws.Cells.Item(46, 5).Value = "  +2.32%  "

# Row 47
$ws.Cells.Item(47, 4).Value = "55.77"
$ws.Cells.Item(47, 5).Value = "  +1.71%  "

# Row 48
$ws.Cells.Item(48, 4).Value = "0.05279"

# Row 49
$ws.Cells.Item(49, 5).Value = "  +0.94%  "

# Row 50
$ws.Cells.Item(50, 2).Value = "Decentraland"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Cells.Item(50, 4).Value = "0.3494"
$ws.Cells.Item(50, 5).Value = "  +2.64%  "

# Row 51
$ws.Cells.Item(51, 2).Value = "EnergySwap"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(51, 4).Value = "7.624"
$ws.Cells.Item(51, 5).Value = "  +0.73%  "
